# Sprint 4 Backlog - Burndown: update actuals for Matthew's tasks and
# refresh the view's selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Week 1 actual hours logged by Matthew, and who completed it.
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "Matthew"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0

# Row 5
$ws.Range("E5").Value = 1.5
$ws.Range("F5").Value = "Matthew"
$ws.Range("H5").Value = 1.5
$ws.Range("I5").Value = 0

# Row 11
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = "Matthew"
$ws.Range("H11").Value = 0.5
$ws.Range("I11").Value = 0

# Row 15: only the "amount remaining" figure was filled in.
$ws.Range("H15").Value = 1

# Update the saved view state (scroll position and active selection).
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E16").Select()
